$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets("ALC")
# Row 17 (ALC)
$ws_ALC.Range("H17").Value = 26233.797
$ws_ALC.Range("J17").Value = 26233.797
$ws_ALC.Range("L17").Value = 78701.391
$ws_ALC.Range("N17").Value = -79037.391

# Row 58 (ALC)
$ws_ALC.Range("H58").Value = 18179.635
$ws_ALC.Range("J58").Value = 21161.797
$ws_ALC.Range("L58").Value = 63485.391
$ws_ALC.Range("N58").Value = -63785.391

# Row 70 (ALC)
$ws_ALC.Range("H70").Value = 1670.4166
$ws_ALC.Range("I70").Value = 969.381
$ws_ALC.Range("J70").Value = 2651.8667
$ws_ALC.Range("K70").Value = 2908.143
$ws_ALC.Range("L70").Value = 7955.6001
$ws_ALC.Range("M70").Value = -2638.143
$ws_ALC.Range("N70").Value = -8495.6001

# Row 73 (ALC)
$ws_ALC.Range("H73").Value = 1670.4166
$ws_ALC.Range("I73").Value = 969.381
$ws_ALC.Range("J73").Value = 2651.8667
$ws_ALC.Range("K73").Value = 2908.143
$ws_ALC.Range("L73").Value = 7955.6001
$ws_ALC.Range("M73").Value = -1972.143
$ws_ALC.Range("N73").Value = -9827.6001

# Row 82 (ALC)
$ws_ALC.Range("H82").Value = 2786.3333

# Row 85 (ALC)
$ws_ALC.Range("H85").Value = 2786.3333

# Row 96 (ALC)
$ws_ALC.Range("H96").Value = 649.75
$ws_ALC.Range("I96").Value = 411
$ws_ALC.Range("K96").Value = 1233
$ws_ALC.Range("M96").Value = 140

# Row 115 (ALC)
$ws_ALC.Range("H115").Value = 952.8570999999999

# Row 121 (ALC)
$ws_ALC.Range("H121").Value = 1110.5454
$ws_ALC.Range("J121").Value = 1571.3334
$ws_ALC.Range("L121").Value = 4714.0002
$ws_ALC.Range("N121").Value = -8208.0002

# Row 138 (ALC)
$ws_ALC.Range("H138").Value = 2340.2659
$ws_ALC.Range("I138").Value = 850.9056399999999
$ws_ALC.Range("J138").Value = 5376.269
$ws_ALC.Range("K138").Value = 2552.71692
$ws_ALC.Range("L138").Value = 16128.807
$ws_ALC.Range("M138").Value = 2587.28308
$ws_ALC.Range("N138").Value = -26408.807

$ws_ARM = $wb.Worksheets("ARM")
# Row 32 (ARM)
$ws_ARM.Range("H32").Value = 926.8
$ws_ARM.Range("I32").Value = 925.30615
$ws_ARM.Range("J32").Value = 1000
$ws_ARM.Range("K32").Value = 925.30615
$ws_ARM.Range("L32").Value = 1000
$ws_ARM.Range("M32").Value = -638.30615
$ws_ARM.Range("N32").Value = -1574

# Row 61 (ARM)
$ws_ARM.Range("H61").Value = 1565.25
$ws_ARM.Range("I61").Value = 864.5357
$ws_ARM.Range("J61").Value = 4017.75
$ws_ARM.Range("K61").Value = 864.5357
$ws_ARM.Range("L61").Value = 4017.75
$ws_ARM.Range("M61").Value = -652.5357
$ws_ARM.Range("N61").Value = -4441.75

# Row 74 (ARM)
$ws_ARM.Range("H74").Value = 1402.2222
$ws_ARM.Range("I74").Value = 945.7143
$ws_ARM.Range("J74").Value = 3000
$ws_ARM.Range("K74").Value = 945.7143
$ws_ARM.Range("L74").Value = 3000
$ws_ARM.Range("M74").Value = -71.71429999999998
$ws_ARM.Range("N74").Value = -4748

# Row 77 (ARM)
$ws_ARM.Range("H77").Value = 1402.2222
$ws_ARM.Range("I77").Value = 945.7143
$ws_ARM.Range("J77").Value = 3000
$ws_ARM.Range("K77").Value = 4728.5715
$ws_ARM.Range("L77").Value = 15000
$ws_ARM.Range("M77").Value = -360.5715
$ws_ARM.Range("N77").Value = -23736

# Row 111 (ARM)
$ws_ARM.Range("H111").Value = 40000
$ws_ARM.Range("J111").Value = 40000
$ws_ARM.Range("L111").Value = 40000

# Row 136 (ARM)
$ws_ARM.Range("H136").Value = 1565.25
$ws_ARM.Range("I136").Value = 864.5357
$ws_ARM.Range("J136").Value = 4017.75
$ws_ARM.Range("K136").Value = 2593.6071
$ws_ARM.Range("L136").Value = 12053.25
$ws_ARM.Range("M136").Value = -43.60710000000017
$ws_ARM.Range("N136").Value = -17153.25

$ws_BSM = $wb.Worksheets("BSM")
# Row 94 (BSM)
$ws_BSM.Range("H94").Value = 936.2632
$ws_BSM.Range("I94").Value = 843.38464
$ws_BSM.Range("J94").Value = 1137.5
$ws_BSM.Range("K94").Value = 843.38464
$ws_BSM.Range("L94").Value = 1137.5
$ws_BSM.Range("M94").Value = -392.38464
$ws_BSM.Range("N94").Value = -2039.5

# Row 105 (BSM)
$ws_BSM.Range("H105").Value = 1511.56
$ws_BSM.Range("I105").Value = 1390.5294
$ws_BSM.Range("J105").Value = 1768.75
$ws_BSM.Range("K105").Value = 1390.5294
$ws_BSM.Range("L105").Value = 1768.75
$ws_BSM.Range("M105").Value = 356.4706000000001
$ws_BSM.Range("N105").Value = -5262.75

$ws_CRP = $wb.Worksheets("CRP")
# Row 31 (CRP)
$ws_CRP.Range("H31").Value = 3652.8572
$ws_CRP.Range("I31").Value = 2502.25
$ws_CRP.Range("J31").Value = 4621.7896
$ws_CRP.Range("K31").Value = 2502.25
$ws_CRP.Range("L31").Value = 4621.7896
$ws_CRP.Range("M31").Value = -2207.25
$ws_CRP.Range("N31").Value = -5211.7896

# Row 34 (CRP)
$ws_CRP.Range("H34").Value = 3652.8572
$ws_CRP.Range("I34").Value = 2502.25
$ws_CRP.Range("J34").Value = 4621.7896
$ws_CRP.Range("K34").Value = 2502.25
$ws_CRP.Range("L34").Value = 4621.7896
$ws_CRP.Range("M34").Value = -2300.25
$ws_CRP.Range("N34").Value = -5025.7896

# Row 58 (CRP)
$ws_CRP.Range("H58").Value = 8335566
$ws_CRP.Range("I58").Value = 1314.8541
$ws_CRP.Range("J58").Value = 41672572
$ws_CRP.Range("K58").Value = 1314.8541
$ws_CRP.Range("L58").Value = 41672572
$ws_CRP.Range("M58").Value = -1111.8541
$ws_CRP.Range("N58").Value = -41672978

# Row 100 (CRP)
$ws_CRP.Range("H100").Value = 0
$ws_CRP.Range("J100").Value = 0
$ws_CRP.Range("L100").Value = 0
$ws_CRP.Range("N100").ClearContents()

# Row 105 (CRP)
$ws_CRP.Range("H105").Value = 2976.25
$ws_CRP.Range("I105").Value = 3015
$ws_CRP.Range("J105").Value = 2937.5
$ws_CRP.Range("K105").Value = 3015
$ws_CRP.Range("L105").Value = 2937.5
$ws_CRP.Range("M105").Value = -1268
$ws_CRP.Range("N105").Value = -6431.5

# Row 132 (CRP)
$ws_CRP.Range("H132").Value = 1787.196
$ws_CRP.Range("I132").Value = 1430.7142
$ws_CRP.Range("K132").Value = 4292.142599999999
$ws_CRP.Range("M132").Value = -1762.142599999999

# Row 134 (CRP)
$ws_CRP.Range("H134").Value = 1438.2642
$ws_CRP.Range("I134").Value = 852.3555
$ws_CRP.Range("J134").Value = 4734
$ws_CRP.Range("K134").Value = 2557.0665
$ws_CRP.Range("L134").Value = 14202
$ws_CRP.Range("M134").Value = -22.06649999999991
$ws_CRP.Range("N134").Value = -19272

# Row 136 (CRP)
$ws_CRP.Range("H136").Value = 8335566
$ws_CRP.Range("I136").Value = 1314.8541
$ws_CRP.Range("J136").Value = 41672572
$ws_CRP.Range("K136").Value = 3944.5623
$ws_CRP.Range("L136").Value = 125017716
$ws_CRP.Range("M136").Value = -1394.5623
$ws_CRP.Range("N136").Value = -125022816

$ws_CUL = $wb.Worksheets("CUL")
# Row 57 (CUL)
$ws_CUL.Range("H57").Value = 3536
$ws_CUL.Range("I57").Value = 1780
$ws_CUL.Range("K57").Value = 5340
$ws_CUL.Range("M57").Value = -4781

# Row 87 (CUL)
$ws_CUL.Range("H87").Value = 5969.9565
$ws_CUL.Range("I87").Value = 2167.2666
$ws_CUL.Range("J87").Value = 13100
$ws_CUL.Range("K87").Value = 6501.7998
$ws_CUL.Range("L87").Value = 39300
$ws_CUL.Range("M87").Value = -5253.7998
$ws_CUL.Range("N87").Value = -41796

# Row 90 (CUL)
$ws_CUL.Range("H90").Value = 5969.9565
$ws_CUL.Range("I90").Value = 2167.2666
$ws_CUL.Range("J90").Value = 13100
$ws_CUL.Range("K90").Value = 19505.3994
$ws_CUL.Range("L90").Value = 117900
$ws_CUL.Range("M90").Value = -13265.3994
$ws_CUL.Range("N90").Value = -130380

# Row 120 (CUL)
$ws_CUL.Range("H120").Value = 17784.777
$ws_CUL.Range("I120").Value = 15676.667
$ws_CUL.Range("K120").Value = 47030.001
$ws_CUL.Range("M120").Value = -42192.001

# Row 134 (CUL)
$ws_CUL.Range("H134").Value = 1894.4117
$ws_CUL.Range("I134").Value = 1019.6667
$ws_CUL.Range("J134").Value = 3993.8
$ws_CUL.Range("K134").Value = 3059.0001
$ws_CUL.Range("L134").Value = 11981.4
$ws_CUL.Range("M134").Value = 2010.9999
$ws_CUL.Range("N134").Value = -22121.4

# Row 137 (CUL)
$ws_CUL.Range("H137").Value = 2733.7812
$ws_CUL.Range("I137").Value = 1925.4667
$ws_CUL.Range("J137").Value = 3447
$ws_CUL.Range("K137").Value = 5776.4001
$ws_CUL.Range("L137").Value = 10341
$ws_CUL.Range("M137").Value = -676.4000999999998
$ws_CUL.Range("N137").Value = -20541

# Row 138 (CUL)
$ws_CUL.Range("H138").Value = 2157.4167
$ws_CUL.Range("I138").Value = 988.4286
$ws_CUL.Range("J138").Value = 3794
$ws_CUL.Range("K138").Value = 2965.2858
$ws_CUL.Range("L138").Value = 11382
$ws_CUL.Range("M138").Value = 2174.7142
$ws_CUL.Range("N138").Value = -21662

# Row 139 (CUL)
$ws_CUL.Range("H139").Value = 9263997
$ws_CUL.Range("I139").Value = 22731450
$ws_CUL.Range("J139").Value = 5122.8125
$ws_CUL.Range("K139").Value = 68194350
$ws_CUL.Range("L139").Value = 15368.4375
$ws_CUL.Range("M139").Value = -68189210
$ws_CUL.Range("N139").Value = -25648.4375

$ws_GSM = $wb.Worksheets("GSM")
# Row 133 (GSM)
$ws_GSM.Range("H133").Value = 35780
$ws_GSM.Range("J133").Value = 35780
$ws_GSM.Range("L133").Value = 35780

$ws_LTW = $wb.Worksheets("LTW")
# Row 46 (LTW)
$ws_LTW.Range("H46").Value = 2104.0908
$ws_LTW.Range("I46").Value = 500
$ws_LTW.Range("J46").Value = 2264.5
$ws_LTW.Range("K46").Value = 500
$ws_LTW.Range("L46").Value = 2264.5
$ws_LTW.Range("M46").Value = -312
$ws_LTW.Range("N46").Value = -2640.5

# Row 132 (LTW)
$ws_LTW.Range("H132").Value = 2131.425
$ws_LTW.Range("I132").Value = 1266.3549
$ws_LTW.Range("K132").Value = 3799.0647
$ws_LTW.Range("M132").Value = -1269.0647

# Row 136 (LTW)
$ws_LTW.Range("H136").Value = 2487.186
$ws_LTW.Range("I136").Value = 1377
$ws_LTW.Range("J136").Value = 9333.333000000001
$ws_LTW.Range("K136").Value = 4131
$ws_LTW.Range("L136").Value = 27999.999
$ws_LTW.Range("M136").Value = -1581
$ws_LTW.Range("N136").Value = -33099.999

$ws_WVR = $wb.Worksheets("WVR")
# Row 132 (WVR)
$ws_WVR.Range("H132").Value = 11576.704
$ws_WVR.Range("I132").Value = 2554.9143
$ws_WVR.Range("J132").Value = 28195.79
$ws_WVR.Range("K132").Value = 7664.742899999999
$ws_WVR.Range("L132").Value = 84587.37
$ws_WVR.Range("M132").Value = -5134.742899999999
$ws_WVR.Range("N132").Value = -89647.37
